$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.91087007522583
$ws.Range("B1").Value = 2.118667125701904
$ws.Range("C1").Value = 2.082835912704468
$ws.Range("D1").Value = 2.492074489593506
$ws.Range("E1").Value = 2.220349073410034
